# Daily attendance processing - 2025-10-23 23:41:18
#
# The "Recorded By" column (G) lists the users/accounts that touched a
# given attendance row. For this pass, every multi-name entry in column G
# has the order of its comma-separated names reversed, except any entry
# that already reads "admin@admin.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 7).End(-4121).Row  # xlDown

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }
    if ($val -eq "admin@admin.com, System") { continue }

    $parts = $val -split ", "
    $count = $parts.Length

    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = [string]::Join(", ", $reversed)
}
